$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1296.8649
$ws.Range("I15").Value = 1296.8649
$ws.Range("K15").Value = 3890.5947
$ws.Range("M15").Value = -3721.5947
$ws.Range("H33").Value = 430.4
$ws.Range("I33").Value = 401.18518
$ws.Range("K33").Value = 401.18518
$ws.Range("M33").Value = -172.18518
$ws.Range("H137").Value = 1637.9166
$ws.Range("I137").Value = 1265.6
$ws.Range("K137").Value = 3796.8
$ws.Range("M137").Value = -1246.8
$ws.Range("H141").Value = 6461.7036
$ws.Range("I141").Value = 7044.45
$ws.Range("J141").Value = 4796.7144
$ws.Range("K141").Value = 21133.35
$ws.Range("L141").Value = 14390.1432
$ws.Range("M141").Value = -15953.35
$ws.Range("N141").Value = -24750.1432

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3724.78
$ws.Range("I32").Value = 3398.6042
$ws.Range("J32").Value = 11553
$ws.Range("K32").Value = 3398.6042
$ws.Range("L32").Value = 11553
$ws.Range("M32").Value = -3111.6042
$ws.Range("N32").Value = -12127
$ws.Range("H45").Value = 4093.9355
$ws.Range("I45").Value = 2560
$ws.Range("K45").Value = 2560
$ws.Range("M45").Value = -2183
$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()
$ws.Range("H61").Value = 4032.3076
$ws.Range("I61").Value = 4940.1816
$ws.Range("K61").Value = 4940.1816
$ws.Range("M61").Value = -4728.1816
$ws.Range("H74").Value = 2508.182
$ws.Range("I74").Value = 2508.182
$ws.Range("K74").Value = 2508.182
$ws.Range("M74").Value = -1634.182
$ws.Range("H77").Value = 2508.182
$ws.Range("I77").Value = 2508.182
$ws.Range("K77").Value = 12540.91
$ws.Range("M77").Value = -8172.91
$ws.Range("H124").Value = 49866.332
$ws.Range("J124").Value = 49866.332
$ws.Range("L124").Value = 49866.332
$ws.Range("N124").Value = -59686.332
$ws.Range("H125").Value = 109354.2
$ws.Range("J125").Value = 109354.2
$ws.Range("L125").Value = 109354.2
$ws.Range("N125").Value = -119194.2
$ws.Range("H136").Value = 4032.3076
$ws.Range("I136").Value = 4940.1816
$ws.Range("K136").Value = 14820.5448
$ws.Range("M136").Value = -12270.5448
$ws.Range("H139").Value = 137399.2
$ws.Range("J139").Value = 137399.2
$ws.Range("L139").Value = 137399.2
$ws.Range("N139").Value = -147679.2

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 150559.86
$ws.Range("J132").Value = 150559.86
$ws.Range("L132").Value = 150559.86
$ws.Range("N132").Value = -160679.86
$ws.Range("H134").Value = 3028.111
$ws.Range("I134").Value = 2712.25
$ws.Range("K134").Value = 8136.75
$ws.Range("M134").Value = -5601.75
$ws.Range("H140").Value = 40385
$ws.Range("J140").Value = 40385
$ws.Range("L140").Value = 40385
$ws.Range("N140").Value = -50745

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1857.8966
$ws.Range("I31").Value = 1503.4166
$ws.Range("J31").Value = 3559.4
$ws.Range("K31").Value = 1503.4166
$ws.Range("L31").Value = 3559.4
$ws.Range("M31").Value = -1208.4166
$ws.Range("N31").Value = -4149.4
$ws.Range("H34").Value = 1857.8966
$ws.Range("I34").Value = 1503.4166
$ws.Range("J34").Value = 3559.4
$ws.Range("K34").Value = 1503.4166
$ws.Range("L34").Value = 3559.4
$ws.Range("M34").Value = -1301.4166
$ws.Range("N34").Value = -3963.4
$ws.Range("H99").Value = 16277.044
$ws.Range("I99").Value = 12987
$ws.Range("J99").Value = 18392.072
$ws.Range("K99").Value = 12987
$ws.Range("L99").Value = 18392.072
$ws.Range("M99").Value = -11489
$ws.Range("N99").Value = -21388.072
$ws.Range("H126").Value = 16277.044
$ws.Range("I126").Value = 12987
$ws.Range("J126").Value = 18392.072
$ws.Range("K126").Value = 38961
$ws.Range("L126").Value = 55176.216
$ws.Range("M126").Value = -36491
$ws.Range("N126").Value = -60116.216
$ws.Range("H134").Value = 2994.5588
$ws.Range("I134").Value = 2820.2334
$ws.Range("K134").Value = 8460.700199999999
$ws.Range("M134").Value = -5925.700199999999

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 806698.2
$ws.Range("I9").Value = 1110011.4
$ws.Range("J9").Value = 200071.8
$ws.Range("K9").Value = 3330034.2
$ws.Range("L9").Value = 600215.3999999999
$ws.Range("M9").Value = -3329810.2
$ws.Range("N9").Value = -600663.3999999999
$ws.Range("H132").Value = 3101.7273
$ws.Range("I132").Value = 2023.25
$ws.Range("J132").Value = 3718
$ws.Range("K132").Value = 18209.25
$ws.Range("L132").Value = 33462
$ws.Range("M132").Value = -15679.25
$ws.Range("N132").Value = -38522
$ws.Range("H136").Value = 1854.3334
$ws.Range("I136").Value = 1854.3334
$ws.Range("K136").Value = 5563.0002
$ws.Range("M136").Value = -463.0002000000004

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 14750
$ws.Range("I43").Value = 4500
$ws.Range("J43").Value = 25000
$ws.Range("K43").Value = 4500
$ws.Range("L43").Value = 25000
$ws.Range("M43").Value = -4349
$ws.Range("N43").Value = -25302
$ws.Range("H123").Value = 2500
$ws.Range("J123").Value = 2500
$ws.Range("L123").Value = 2500
$ws.Range("N123").Value = -7400
$ws.Range("H132").Value = 2531.303
$ws.Range("I132").Value = 1618.6666
$ws.Range("K132").Value = 4855.9998
$ws.Range("M132").Value = -2325.9998
$ws.Range("H135").Value = 145177
$ws.Range("J135").Value = 145177
$ws.Range("L135").Value = 145177
$ws.Range("N135").Value = -155317

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6154.7
$ws.Range("I136").Value = 6720.893
$ws.Range("J136").Value = 4833.5835
$ws.Range("K136").Value = 20162.679
$ws.Range("L136").Value = 14500.7505
$ws.Range("M136").Value = -17612.679
$ws.Range("N136").Value = -19600.7505
$ws.Range("H139").Value = 133820.5
$ws.Range("J139").Value = 128760.664
$ws.Range("L139").Value = 128760.664
$ws.Range("N139").Value = -139040.664

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 8815.736999999999
$ws.Range("I122").Value = 6107.2144
$ws.Range("K122").Value = 18321.6432
$ws.Range("M122").Value = -15871.6432
$ws.Range("H139").Value = 66299.14
$ws.Range("J139").Value = 66299.14
$ws.Range("L139").Value = 66299.14
$ws.Range("N139").Value = -76579.14
